# Apply the data refresh captured in the commit:
# "Update gh-pages to output generated at 456a3b4"
#
# The workbook lists Guangzhou comic/anime-convention events across four
# sheets ("展览" Exhibitions, "演出" Performances, "本地生活" Local Life,
# "全部类型" All Types — the last one aggregates rows from the first
# three). Column F holds a view/click counter and column G holds the
# ticket price (or a status string such as "暂时售罄" / temporarily sold
# out when the price isn't available). This script pushes the refreshed
# counters/status pulled from the live site into the matching cells on
# every sheet.

$wb = $excel.ActiveWorkbook

function Set-Cell($ws, $addr, $value) {
    $ws.Range($addr).Value = $value
}

# --- 展览 (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
Set-Cell $wsExpo "F3"  1346
Set-Cell $wsExpo "F4"  13214
Set-Cell $wsExpo "F11" 47
Set-Cell $wsExpo "F13" 19091
Set-Cell $wsExpo "G13" "暂时售罄"
Set-Cell $wsExpo "F19" 312
Set-Cell $wsExpo "F21" 137
Set-Cell $wsExpo "F25" 4
Set-Cell $wsExpo "F26" 1341

# --- 演出 (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
Set-Cell $wsShow "F5"  186
Set-Cell $wsShow "F8"  85
Set-Cell $wsShow "F9"  85
Set-Cell $wsShow "F10" 379

# --- 本地生活 (Local Life) ---
$wsLocal = $wb.Worksheets.Item("本地生活")
Set-Cell $wsLocal "F3" 4379

# --- 全部类型 (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
Set-Cell $wsAll "F5"  1346
Set-Cell $wsAll "F6"  13214
Set-Cell $wsAll "F9"  4379
Set-Cell $wsAll "F15" 47
Set-Cell $wsAll "F17" 19091
Set-Cell $wsAll "G17" "暂时售罄"
Set-Cell $wsAll "F21" 186
Set-Cell $wsAll "F22" 186
Set-Cell $wsAll "F26" 85
Set-Cell $wsAll "F27" 85
Set-Cell $wsAll "F28" 379
Set-Cell $wsAll "F32" 312
Set-Cell $wsAll "F34" 137
Set-Cell $wsAll "F41" 4
Set-Cell $wsAll "F42" 1341
